$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20121215D")

# Row 1 headers
$ws.Range("C1").Value = "Read only string in parser and symbol"
$ws.Range("D1").Value = "Remove unused symbol insertion when parameter is anonymous"
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""

# Column D data values (rows 2-11)
$ws.Range("D2").Value = 4765
$ws.Range("D3").Value = 4747
$ws.Range("D4").Value = 4771
$ws.Range("D5").Value = 4754
$ws.Range("D6").Value = 4751
$ws.Range("D7").Value = 4812
$ws.Range("D8").Value = 4771
$ws.Range("D9").Value = 4795
$ws.Range("D10").Value = 4736
$ws.Range("D11").Value = 4780

# Update selection to D15
$ws.Range("D15").Select()
